$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Source URL text updates (shared strings) ---

# Denmark - tests performed (row 22)
$ws.Range("D22").Value = "https://files.ssi.dk/Data-Epidemiologiske-Rapport-10082020-lz41"

# Iran - tests performed (row 40)
$ws.Range("D40").Value = "http://irangov.ir/detail/344872"

# Spain - tests performed (row 86)
$ws.Range("D86").Value = "https://www.mscbs.gob.es/profesionales/saludPublica/ccayes/alertasActual/nCov-China/documentos/COVID-19_pruebas_diagnosticas_06_08_2020.pdf"

# --- Numeric data updates ---

# Row 18 - Croatia - people tested
$ws.Range("C18").Value = 44053
$ws.Range("G18").Value = 152
$ws.Range("H18").Value = 129379
$ws.Range("I18").Value = 31.515
$ws.Range("J18").Value = 1247
$ws.Range("K18").Value = 0.304
$ws.Range("L18").Value = 1042
$ws.Range("M18").Value = 0.254
$ws.Range("O18").Value = 21.203

# Row 22 - Denmark - tests performed
$ws.Range("C22").Value = 44052
$ws.Range("G22").Value = 183
$ws.Range("H22").Value = 1721479
$ws.Range("I22").Value = 297.206
$ws.Range("J22").Value = 3826
$ws.Range("K22").Value = 0.661
$ws.Range("L22").Value = 21501
$ws.Range("M22").Value = 3.712
$ws.Range("O22").Value = 230.485

# Row 36 - Iceland - tests performed
$ws.Range("C36").Value = 44052
$ws.Range("G36").Value = 165
$ws.Range("H36").Value = 78467
$ws.Range("I36").Value = 229.94
$ws.Range("J36").Value = 232
$ws.Range("K36").Value = 0.68
$ws.Range("L36").Value = 644
$ws.Range("M36").Value = 1.887
$ws.Range("N36").Value = 0.014
$ws.Range("O36").Value = 72.71

# Row 40 - Iran - tests performed
$ws.Range("C40").Value = 44053
$ws.Range("G40").Value = 111
$ws.Range("H40").Value = 2711817
$ws.Range("I40").Value = 32.286
$ws.Range("J40").Value = 25319
$ws.Range("K40").Value = 0.301
$ws.Range("L40").Value = 25308
$ws.Range("M40").Value = 0.301
$ws.Range("N40").Value = 0.098
$ws.Range("O40").Value = 10.255

# Row 62 - Nigeria - samples tested
$ws.Range("C62").Value = 44053
$ws.Range("G62").Value = 95
$ws.Range("H62").Value = 319851
$ws.Range("I62").Value = 1.552
$ws.Range("J62").Value = 2355
$ws.Range("K62").Value = 0.011
$ws.Range("L62").Value = 4388
$ws.Range("N62").Value = 0.089
$ws.Range("O62").Value = 11.227

# Row 63 - Norway - people tested
$ws.Range("C63").Value = 44051
$ws.Range("G63").Value = 167
$ws.Range("H63").Value = 474531
$ws.Range("I63").Value = 87.532
$ws.Range("J63").Value = 1990
$ws.Range("K63").Value = 0.367
$ws.Range("L63").Value = 5463
$ws.Range("M63").Value = 1.008
$ws.Range("N63").Value = 0.007
$ws.Range("O63").Value = 147.081

# Row 86 - Spain - tests performed
$ws.Range("C86").Value = 44049
$ws.Range("G86").Value = 17
$ws.Range("H86").Value = 4983935
$ws.Range("I86").Value = 106.597
$ws.Range("L86").Value = 47349
$ws.Range("M86").Value = 1.013
$ws.Range("N86").Value = 0.074
$ws.Range("O86").Value = 13.57
